# edit.ps1 - Applies the "cambios de agosto, puntos fe de ratas e historico" update
# to a69_f34_dUPPachuca.xlsx ("Reporte de Formatos" - Inventario de bienes inmuebles)
#
# Main changes:
#   * Update reporting period from 2021 (3er/4o trimestre) to 2022 (1er/2o trimestre)
#   * Update validation / update dates
#   * Replace the explanatory note text in column AI ("Nota")
#   * Re-wrap the long description header (G3:I3) and tidy its font/border
#   * Extend the Hidden_* list-validation ranges from row 90 to row 201
#   * Narrow column AI and move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# ---------------------------------------------------------------------------
# 1) Row 8 - the single data row of the report
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 2022
$ws.Range("B8").Value = 44562   # 2022-01-01
$ws.Range("C8").Value = 44742   # 2022-06-30

# AF8 keeps its text ("Secretaría Administrativa (UPP)") but gets a plain
# black Calibri font (no longer inherits the old left-aligned style)
$ws.Range("AF8").Font.Name = "Calibri"
$ws.Range("AF8").Font.Color = RGB(0, 0, 0)
$ws.Range("AF8").HorizontalAlignment = 1

$ws.Range("AG8").Value = 44753   # 2022-07-11
$ws.Range("AH8").Value = 44753   # 2022-07-11

# The note text is replaced completely
$ws.Range("AI8").Value = "La Universidad Politécnica de Pachuca, no cuenta con inventario de bienes inmuebles. "

# ---------------------------------------------------------------------------
# 2) Header area (row 3) - wrap the long description text and tidy borders
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 46.5

$ws.Range("H3:I3").WrapText = $true

$ws.Range("G3").Borders.LineStyle = -4142   # xlLineStyleNone - clear all edges
$ws.Range("G3").Borders.Item(7).LineStyle = 1   # xlEdgeLeft = thin
$ws.Range("G3").WrapText = $true

# ---------------------------------------------------------------------------
# 3) Column AI ("Nota") gets narrower now that the note is shorter
# ---------------------------------------------------------------------------
$ws.Columns.Item(35).ColumnWidth = 47.02

# ---------------------------------------------------------------------------
# 4) Extend the six list-validations from row 90 to row 201
# ---------------------------------------------------------------------------
$ws.Range("F8:F90").Validation.Delete()
$ws.Range("F8:F201").Validation.Add(3, 1, 1, "Hidden_15")
$ws.Range("F8:F201").Validation.ShowInput = $false

$ws.Range("J8:J90").Validation.Delete()
$ws.Range("J8:J201").Validation.Add(3, 1, 1, "Hidden_29")
$ws.Range("J8:J201").Validation.ShowInput = $false

$ws.Range("Q8:Q90").Validation.Delete()
$ws.Range("Q8:Q201").Validation.Add(3, 1, 1, "Hidden_316")
$ws.Range("Q8:Q201").Validation.ShowInput = $false

$ws.Range("W8:W90").Validation.Delete()
$ws.Range("W8:W201").Validation.Add(3, 1, 1, "Hidden_422")
$ws.Range("W8:W201").Validation.ShowInput = $false

$ws.Range("X8:X90").Validation.Delete()
$ws.Range("X8:X201").Validation.Add(3, 1, 1, "Hidden_523")
$ws.Range("X8:X201").Validation.ShowInput = $false

$ws.Range("Y8:Y90").Validation.Delete()
$ws.Range("Y8:Y201").Validation.Add(3, 1, 1, "Hidden_624")
$ws.Range("Y8:Y201").Validation.ShowInput = $false

# ---------------------------------------------------------------------------
# 5) Move the active selection (the user left the cursor on A9)
# ---------------------------------------------------------------------------
$ws.Range("A9").Select()
